# Updates the "Price" column (D) values for a set of rows in the crypto
# symbol list, as published by the "Updated symbol list" GitHub Action.
#
# The values are written as literal text (matching the source data, which
# stores prices as inline strings so exact formatting - e.g. trailing
# zeros - is preserved) rather than as numbers, which would otherwise get
# reformatted/rounded by Excel's numeric type coercion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2";  Value = "266.37" }
    @{ Cell = "D3";  Value = "21.32" }
    @{ Cell = "D4";  Value = "6.116" }
    @{ Cell = "D5";  Value = "0.06105" }
    @{ Cell = "D6";  Value = "3.566" }
    @{ Cell = "D7";  Value = "6.485" }
    @{ Cell = "D8";  Value = "1.355" }
    @{ Cell = "D9";  Value = "0.8208" }
    @{ Cell = "D10"; Value = "0.01338" }
    @{ Cell = "D11"; Value = "0.1593" }
    @{ Cell = "D12"; Value = "0.08023" }
    @{ Cell = "D14"; Value = "0.03218" }
    @{ Cell = "D16"; Value = "3.735" }
    @{ Cell = "D17"; Value = "0.001623" }
    @{ Cell = "D18"; Value = "0.04656" }
    @{ Cell = "D19"; Value = "0.006312" }
    @{ Cell = "D20"; Value = "0.006147" }
    @{ Cell = "D21"; Value = "0.001069" }
    @{ Cell = "D23"; Value = "3.726" }
    @{ Cell = "D25"; Value = "0.3313" }
    @{ Cell = "D28"; Value = "0.0002713" }
    @{ Cell = "D40"; Value = "0.04601" }
    @{ Cell = "D41"; Value = "0.006990" }
    @{ Cell = "D42"; Value = "0.004000" }
    @{ Cell = "D44"; Value = "0.01057" }
    @{ Cell = "D45"; Value = "0.00005949" }
    @{ Cell = "D46"; Value = "0.0009902" }
    @{ Cell = "D48"; Value = "0.8026" }
    @{ Cell = "D49"; Value = "0.001126" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text formatting so the numeric-looking string is not
    # auto-converted into a floating point number (which would drop
    # significant trailing zeros / change precision).
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    # Reset the style back to Normal so no stray cell formatting is left
    # behind - only the textual content should differ from the original.
    $cell.Style = "Normal"
}

$wb.Save()
